$wb = $excel.ActiveWorkbook

# --- sheet3 ("attributes"): insert the missing "ID" row for the
#     rd_bb_contribution entity group, as the first row of that group
#     (row 120), shifting all following rows down by one. ---
$wsAttrs = $wb.Worksheets.Item("attributes")
$wsAttrs.Rows.Item(120).Insert()
$wsAttrs.Cells.Item(120, 1).Value = "ID"
$wsAttrs.Cells.Item(120, 2).Value = "ID"
$wsAttrs.Cells.Item(120, 3).Value = " "
$wsAttrs.Cells.Item(120, 4).Value = "rd_bb_contribution"
# Copy the literal text "true" from an existing ID row (rather than
# assigning the string "true" directly) so Excel keeps it as text
# instead of auto-converting it to a boolean TRUE value.
$wsAttrs.Cells.Item(131, 8).Copy($wsAttrs.Cells.Item(120, 8))
$wsAttrs.Cells.Item(131, 8).Copy($wsAttrs.Cells.Item(120, 10))

# --- sheet8 ("rd_bb_contribution"): add a new "ID" column as the
#     first column of the header row, shifting the existing headers
#     right by one column. ---
$wsContrib = $wb.Worksheets.Item("rd_bb_contribution")
$wsContrib.Columns.Item(1).Insert()
$wsContrib.Cells.Item(1, 1).Value = "ID"
